$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 757.2727
$ws.Range("I2").Value = 260.33334
$ws.Range("J2").Value = 943.625
$ws.Range("K2").Value = 260.33334
$ws.Range("L2").Value = 943.625
$ws.Range("M2").Value = -147.33334
$ws.Range("N2").Value = -1169.625

$ws.Range("H13").Value = 766.3333
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 766.3333
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 766.3333
$ws.Range("N13").Value = -1104.3333
$ws.Range("M13").ClearContents()

$ws.Range("H34").Value = 11099
$ws.Range("I34").Value = 11099
$ws.Range("K34").Value = 11099
$ws.Range("M34").Value = -10896

$ws.Range("H36").Value = 11099
$ws.Range("I36").Value = 11099
$ws.Range("K36").Value = 11099
$ws.Range("M36").Value = -10384

$ws.Range("H52").Value = 20
$ws.Range("I52").Value = 20
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 60
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = 100
$ws.Range("N52").ClearContents()

$ws.Range("H98").Value = 4006
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").ClearContents()

$ws.Range("H122").Value = 4006
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws.Range("H135").Value = 1602.8
$ws.Range("I135").Value = 1336.6154
$ws.Range("K135").Value = 12029.5386
$ws.Range("M135").Value = -9494.5386

$ws.Range("H138").Value = 4665
$ws.Range("J138").Value = 5198
$ws.Range("L138").Value = 15594
$ws.Range("N138").Value = -25874

$ws.Range("H141").Value = 5111.8125
$ws.Range("I141").Value = 5596.3335
$ws.Range("K141").Value = 16789.0005
$ws.Range("M141").Value = -11609.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 2750
$ws.Range("I19").Value = 3000
$ws.Range("K19").Value = 3000
$ws.Range("M19").Value = -2771

$ws.Range("H56").Value = 35000
$ws.Range("J56").Value = 35000
$ws.Range("L56").Value = 35000
$ws.Range("N56").Value = -36484

$ws.Range("H61").Value = 3848.0334
$ws.Range("J61").Value = 9350.166999999999
$ws.Range("L61").Value = 9350.166999999999
$ws.Range("N61").Value = -9774.166999999999

$ws.Range("H122").Value = 502873.1
$ws.Range("I122").Value = 770804.75
$ws.Range("J122").Value = 5285.7144
$ws.Range("K122").Value = 2312414.25
$ws.Range("L122").Value = 15857.1432
$ws.Range("M122").Value = -2309964.25
$ws.Range("N122").Value = -20757.1432

$ws.Range("H132").Value = 18879.262
$ws.Range("I132").Value = 19378.684
$ws.Range("K132").Value = 58136.052
$ws.Range("M132").Value = -55606.052

$ws.Range("H136").Value = 3848.0334
$ws.Range("J136").Value = 9350.166999999999
$ws.Range("L136").Value = 28050.501
$ws.Range("N136").Value = -33150.501

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 1500
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

$ws.Range("H76").Value = 22822.309
$ws.Range("J76").Value = 22822.309
$ws.Range("L76").Value = 22822.309
$ws.Range("N76").Value = -23452.309

$ws.Range("H79").Value = 22822.309
$ws.Range("J79").Value = 22822.309
$ws.Range("L79").Value = 22822.309
$ws.Range("N79").Value = -25006.309

$ws.Range("H115").Value = 75000
$ws.Range("J115").Value = 75000
$ws.Range("L115").Value = 75000
$ws.Range("N115").Value = -78134

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 784.6
$ws.Range("I32").Value = 499.66666
$ws.Range("K32").Value = 499.66666
$ws.Range("M32").Value = -183.66666

$ws.Range("H39").Value = 3499.5
$ws.Range("I39").Value = 3499.5
$ws.Range("K39").Value = 3499.5
$ws.Range("M39").Value = -3108.5

$ws.Range("H45").Value = 6400
$ws.Range("I45").Value = 6400
$ws.Range("K45").Value = 6400
$ws.Range("M45").Value = -5807

$ws.Range("H49").Value = 3499.5
$ws.Range("I49").Value = 3499.5
$ws.Range("K49").Value = 3499.5
$ws.Range("M49").Value = -3317.5

$ws.Range("H58").Value = 2834.6956
$ws.Range("I58").Value = 1960.7778
$ws.Range("J58").Value = 5980.8
$ws.Range("K58").Value = 1960.7778
$ws.Range("L58").Value = 5980.8
$ws.Range("M58").Value = -1757.7778
$ws.Range("N58").Value = -6386.8

$ws.Range("H59").Value = 21582.2
$ws.Range("J59").Value = 50000
$ws.Range("L59").Value = 50000
$ws.Range("N59").Value = -52290

$ws.Range("H136").Value = 2834.6956
$ws.Range("I136").Value = 1960.7778
$ws.Range("J136").Value = 5980.8
$ws.Range("K136").Value = 5882.3334
$ws.Range("L136").Value = 17942.4
$ws.Range("M136").Value = -3332.3334
$ws.Range("N136").Value = -23042.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 34000276
$ws.Range("I4").Value = 37400212
$ws.Range("K4").Value = 112200636
$ws.Range("M4").Value = -112200524

$ws.Range("H38").Value = 285.54544
$ws.Range("I38").Value = 75
$ws.Range("J38").Value = 847
$ws.Range("K38").Value = 225
$ws.Range("L38").Value = 2541
$ws.Range("M38").Value = 122
$ws.Range("N38").Value = -3235

$ws.Range("H107").Value = 1195.2
$ws.Range("J107").Value = 1357.125
$ws.Range("L107").Value = 4071.375
$ws.Range("N107").Value = -7911.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 3476.5
$ws.Range("J7").Value = 3476.5
$ws.Range("L7").Value = 3476.5
$ws.Range("N7").Value = -3700.5

$ws.Range("H8").Value = 3476.5
$ws.Range("J8").Value = 3476.5
$ws.Range("L8").Value = 3476.5
$ws.Range("N8").Value = -3754.5

$ws.Range("H62").Value = 39975
$ws.Range("I62").Value = 39950
$ws.Range("K62").Value = 39950
$ws.Range("M62").Value = -39264

$ws.Range("H65").Value = 39975
$ws.Range("I65").Value = 39950
$ws.Range("K65").Value = 119850
$ws.Range("M65").Value = -116418

$ws.Range("H122").Value = 106308.1
$ws.Range("I122").Value = 6798.143
$ws.Range("J122").Value = 338498
$ws.Range("K122").Value = 20394.429
$ws.Range("L122").Value = 1015494
$ws.Range("M122").Value = -17944.429
$ws.Range("N122").Value = -1020394

$ws.Range("H132").Value = 12265.523
$ws.Range("I132").Value = 11976.889
$ws.Range("K132").Value = 35930.667
$ws.Range("M132").Value = -33400.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

$ws.Range("H32").Value = 753324.25
$ws.Range("I32").Value = 753324.25
$ws.Range("K32").Value = 753324.25
$ws.Range("M32").Value = -753007.25

$ws.Range("H61").Value = 14705.348
$ws.Range("I61").Value = 13437
$ws.Range("K61").Value = 13437
$ws.Range("M61").Value = -13235

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H93").Value = 15664.5
$ws.Range("I93").Value = 13797.4
$ws.Range("K93").Value = 13797.4
$ws.Range("M93").Value = -12549.4

$ws.Range("H113").Value = 14705.348
$ws.Range("I113").Value = 13437
$ws.Range("K113").Value = 13437
$ws.Range("M113").Value = -11267

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 60000
$ws.Range("I51").Value = 60000
$ws.Range("K51").Value = 60000
$ws.Range("M51").Value = -59490

$ws.Range("H75").Value = 36697.8
$ws.Range("J75").Value = 36624.75
$ws.Range("L75").Value = 36624.75
$ws.Range("N75").Value = -38496.75

$ws.Range("H78").Value = 36697.8
$ws.Range("J78").Value = 36624.75
$ws.Range("L78").Value = 109874.25
$ws.Range("N78").Value = -119234.25

$ws.Range("H81").Value = 1987.5
$ws.Range("I81").Value = 2028.5714
$ws.Range("K81").Value = 4057.1428
$ws.Range("M81").Value = -2996.1428

$ws.Range("H84").Value = 1987.5
$ws.Range("I84").Value = 2028.5714
$ws.Range("K84").Value = 20285.714
$ws.Range("M84").Value = -14981.714

$ws.Range("H113").Value = 1577.8
$ws.Range("I113").Value = 1496.3334
$ws.Range("J113").Value = 1700
$ws.Range("K113").Value = 4489.0002
$ws.Range("L113").Value = 5100
$ws.Range("M113").Value = -2319.0002
$ws.Range("N113").Value = -9440

$ws.Range("H132").Value = 3958
$ws.Range("I132").Value = 3958
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11874
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9344
$ws.Range("N132").ClearContents()
